$wb = $excel.ActiveWorkbook

# --- "Test Cases" sheet: update Runmode for the last test case to "Y"
# and move the sheet's selection down to the row that was just edited.
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("C6").Value = "Y"
$wsTestCases.Range("B6").Select()

# --- "Customer" sheet: rename the last column header so it's specific to
# the Customer sheet, then move the selection to where work continues.
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Range("L1").Value = "Phone_C"
$wsCustomer.Range("J4").Select()

# --- Finish on "NoviceTester", which becomes the active tab.
$wsNoviceTester = $wb.Worksheets.Item("NoviceTester")
$wsNoviceTester.Activate()
$wsNoviceTester.Range("C11").Select()
